{"js": "// Apply the targeted text edits described by the diff using Word's\n// Office.js find/replace (Range.search + Range.insertText).\nconst replacements = [\n  [\"Ativa\u00e7\u00e3o: 01/01/2021\", \"Ativa\u00e7\u00e3o: 01/01/2024\"],\n  [\"Introduzir os conceitos fundamentais da ci\u00eancia administra\u00e7\u00e3o e de configura\u00e7\u00f5es de uma organiza\u00e7\u00e3o.\", \"Introduzir os conceitos fundamentais de administra\u00e7\u00e3o, de configura\u00e7\u00f5es de uma organiza\u00e7\u00e3o, de marketing e Gest\u00e3o de Pessoas, de forma gen\u00e9rica. A disciplina privilegia a discuss\u00e3o dos fundamentos das diversas abordagens e linhas de pensamento administrativo, sob a \u00f3tica da engenharia.\"],\n  [\"Introduce the fundamental concepts of management science and organization selttings.\", \"To Introduce the fundamental concepts of administration, configurations of an organization, marketing and People Management, in a generic way. The subject privileges the discussion of the fundamentals of the different approaches and lines of administrative thought, from the perspective of engineering.\"],\n  [\"1. \u00c1reas de Atua\u00e7\u00e3o da Administra\u00e7\u00e3o.2. Estrutura organizacional.\", \"Conceitos fundamentais de administra\u00e7\u00e3o e no\u00e7\u00f5es b\u00e1sicas de marketing e Gest\u00e3o de Pessoas.\"],\n  [\"1. Management Practice Areas. 2. Organizational structure\", \"Fundamental concepts of administration and basic notions of marketing and People Management.\"],\n  [\"1. No\u00e7\u00f5es b\u00e1sicas de Marketing, Finan\u00e7as e Recursos Humanos. 2. Diferentes configura\u00e7\u00f5es de organiza\u00e7\u00e3o.\", \"1. Elementos de organiza\u00e7\u00f5es de alto desempenho: aprendizagem organizacional, modelo da compet\u00eancia e capacita\u00e7\u00f5es din\u00e2micas.2. Introdu\u00e7\u00e3o \u00e0 Gest\u00e3o de Pessoas3. No\u00e7\u00f5es b\u00e1sicas de Marketing4. Desenvolvimento de atividade pr\u00e1tica extensionista junto \u00e0 micro e pequenos empreendedores da regi\u00e3o (componente curricular: plano de marketing)5. Visita (viagem did\u00e1tica complementar) a uma empresa para conhecer e entender os diferentes processos organizacionais.\"],\n  [\"1. Basic notions of Marketing, Finance and Human Resources.2. Different organization settings.\", \"1. Elements of high performance organizations: organizational learning, competency model and dynamic capabilities.2. Introduction to People Management3. Marketing basics4. Development of practical extension activities with micro and small entrepreneurs in the region (curricular component: marketing plan)5. Visit (complementary didactic trip) to a company to know and understand the different organizational processes.\"],\n  [\"GITMAN, L. J. - ZUTTER, C. J. Princ\u00edpios de Administra\u00e7\u00e3o Financeira. 14 ed. S\u00e3o Paulo: Perason, 2017.GROPPELLI, A. A.; NIKBAKHT, E. Administra\u00e7\u00e3o Financeira. 3 ed. S\u00e3o Paulo: Saraiva, 2010.MARCOUS\u00c9, I.; SURRIDGE, M.; GILLESPIE, A. Finan\u00e7as. S\u00e3o Paulo: Saraiva, 2013.\", \"\"],\n  [\"MORGAN, G. Imagens da organiza\u00e7\u00e3o. S\u00e3o Paulo, Atlas, 1996.\", \"\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText.substring(0, 60));\n  }\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the targeted text edits described by the diff using Word's\n# Find/Replace (Range.Find.Execute with wdReplaceAll).\n$d = $word.ActiveDocument\n\nfunction Replace-Text($doc, $findText, $replaceText) {\n  $find = $doc.Content.Find\n  $ok = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n  if (-not $ok) {\n    throw \"Find/Replace failed for: $findText\"\n  }\n}\n\nReplace-Text $d \"Ativa\u00e7\u00e3o: 01/01/2021\" \"Ativa\u00e7\u00e3o: 01/01/2024\"\nReplace-Text $d \"Introduzir os conceitos fundamentais da ci\u00eancia administra\u00e7\u00e3o e de configura\u00e7\u00f5es de uma organiza\u00e7\u00e3o.\" \"Introduzir os conceitos fundamentais de administra\u00e7\u00e3o, de configura\u00e7\u00f5es de uma organiza\u00e7\u00e3o, de marketing e Gest\u00e3o de Pessoas, de forma gen\u00e9rica. A disciplina privilegia a discuss\u00e3o dos fundamentos das diversas abordagens e linhas de pensamento administrativo, sob a \u00f3tica da engenharia.\"\nReplace-Text $d \"Introduce the fundamental concepts of management science and organization selttings.\" \"To Introduce the fundamental concepts of administration, configurations of an organization, marketing and People Management, in a generic way. The subject privileges the discussion of the fundamentals of the different approaches and lines of administrative thought, from the perspective of engineering.\"\nReplace-Text $d \"1. \u00c1reas de Atua\u00e7\u00e3o da Administra\u00e7\u00e3o.2. Estrutura organizacional.\" \"Conceitos fundamentais de administra\u00e7\u00e3o e no\u00e7\u00f5es b\u00e1sicas de marketing e Gest\u00e3o de Pessoas.\"\nReplace-Text $d \"1. Management Practice Areas. 2. Organizational structure\" \"Fundamental concepts of administration and basic notions of marketing and People Management.\"\nReplace-Text $d \"1. No\u00e7\u00f5es b\u00e1sicas de Marketing, Finan\u00e7as e Recursos Humanos. 2. Diferentes configura\u00e7\u00f5es de organiza\u00e7\u00e3o.\" \"1. Elementos de organiza\u00e7\u00f5es de alto desempenho: aprendizagem organizacional, modelo da compet\u00eancia e capacita\u00e7\u00f5es din\u00e2micas.2. Introdu\u00e7\u00e3o \u00e0 Gest\u00e3o de Pessoas3. No\u00e7\u00f5es b\u00e1sicas de Marketing4. Desenvolvimento de atividade pr\u00e1tica extensionista junto \u00e0 micro e pequenos empreendedores da regi\u00e3o (componente curricular: plano de marketing)5. Visita (viagem did\u00e1tica complementar) a uma empresa para conhecer e entender os diferentes processos organizacionais.\"\nReplace-Text $d \"1. Basic notions of Marketing, Finance and Human Resources.2. Different organization settings.\" \"1. Elements of high performance organizations: organizational learning, competency model and dynamic capabilities.2. Introduction to People Management3. Marketing basics4. Development of practical extension activities with micro and small entrepreneurs in the region (curricular component: marketing plan)5. Visit (complementary didactic trip) to a company to know and understand the different organizational processes.\"\nReplace-Text $d \"GITMAN, L. J. - ZUTTER, C. J. Princ\u00edpios de Administra\u00e7\u00e3o Financeira. 14 ed. S\u00e3o Paulo: Perason, 2017.GROPPELLI, A. A.; NIKBAKHT, E. Administra\u00e7\u00e3o Financeira. 3 ed. S\u00e3o Paulo: Saraiva, 2010.MARCOUS\u00c9, I.; SURRIDGE, M.; GILLESPIE, A. Finan\u00e7as. S\u00e3o Paulo: Saraiva, 2013.\" \"\"\nReplace-Text $d \"MORGAN, G. Imagens da organiza\u00e7\u00e3o. S\u00e3o Paulo, Atlas, 1996.\" \"\"\n"}
